$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.506.03"
$ws.Range("E2").Value = "  +5.54%  "
$ws.Range("D3").Value = "1.725.57"
$ws.Range("E3").Value = "  +4.79%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.59"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5353"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2667"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06585"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.58"
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.606"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.726.57"
$ws.Range("E13").Value = "  +4.88%  "
$ws.Range("D14").Value = "1.964.30"
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").Value = "0.0₅8266"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.81"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "27.522.09"
$ws.Range("E18").Value = "  +5.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.85"
$ws.Range("E19").Value = "  +12.95%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.727"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.067"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.27"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.757"
$ws.Range("E26").Value = "  +16.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1236"
$ws.Range("E27").Value = "  +4.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.365"
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.49"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05491"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.300"
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.562"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.443"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("E34").Value = "  +7.15%  "
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9624"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.425"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01647"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.899"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.051.64"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8483"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.37"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "1.870.15"
$ws.Range("E45").Value = "  +4.91%  "
$ws.Range("D46").Value = "0.0₈119"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.78"
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("E48").Value = "  +3.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.197"
$ws.Range("E49").Value = "  +3.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  +0.25%  "
